$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
$data = @(
    @(1.02, 1.032430907700835, 1.042042578315674, 1.04220355530835, 1.053943330564776, 1.036131861345853, 1.03756116767077, 1.044820383522397, 1.044980905521385, 1.056687905001926, 1.016489487714868),
    @(1.02, 1.033280417173075, 1.042693878745581, 1.042940879729541, 1.054722026554509, 1.036258365324574, 1.038053476917823, 1.045282853160563, 1.045529206319954, 1.057279834166451, 1.016654405083623),
    @(1.02, 1.033830886814595, 1.043115817585385, 1.043418997515129, 1.055226821826519, 1.036339108716106, 1.038372159457525, 1.04558191930708, 1.045884343985034, 1.057663116884903, 1.016761111502511),
    @(1.02, 1.034062489332263, 1.043293319393685, 1.043620240388966, 1.055439257586884, 1.036372786033922, 1.038506162304154, 1.045707601543121, 1.046033726229317, 1.057824310672104, 1.016805968911253),
    @(1.02, 1.034101387217841, 1.043323129634356, 1.043654044096944, 1.05547493931633, 1.036378424914016, 1.03852866361597, 1.045728701452532, 1.046058812951651, 1.057851379388787, 1.016813500537559),
    @(1.02, 1.033833980773661, 1.043118188909161, 1.043421685583181, 1.055229659542381, 1.036339559764897, 1.038373949899013, 1.045583598858586, 1.045886339716138, 1.057665270523751, 1.016761710897769),
    @(1.02, 1.032717840832299, 1.042262582855372, 1.042452525275738, 1.054206301354345, 1.036174844010602, 1.037727519082762, 1.044976714246408, 1.045166132884579, 1.056887894348424, 1.016545223284965),
    @(1.02, 1.030757110537028, 1.040758839754112, 1.040752635234142, 1.052410204058751, 1.035876100150645, 1.036589446423502, 1.04390597133438, 1.043899786852262, 1.055520162777362, 1.016163719400605),
    @(1.02, 1.029454126901467, 1.039759110612371, 1.039624798336864, 1.051217764249046, 1.035671269701001, 1.035831498194754, 1.043191331874102, 1.043057495729703, 1.054609858345312, 1.015909396361743),
    @(1.02, 1.028890931052661, 1.039326897208558, 1.039137743254492, 1.050702626927263, 1.035581240536832, 1.035503497119342, 1.042881711061414, 1.042693253517065, 1.054216067692173, 1.015799281014891),
    @(1.02, 1.028681887700136, 1.039166457421295, 1.038957027473344, 1.05051146395435, 1.03554759965709, 1.035381693659086, 1.042766678955975, 1.042558030712045, 1.054069854715967, 1.015758381012072),
    @(1.02, 1.028726721272896, 1.039200867584454, 1.038995782622648, 1.050552460792464, 1.035554824780739, 1.035407819501581, 1.042791354847498, 1.042587033135847, 1.054101215224783, 1.015767154114649),
    @(1.02, 1.028873648339695, 1.039313633082649, 1.039122801177689, 1.050686821611942, 1.035578463847835, 1.035493428172241, 1.042872202976219, 1.042682074460298, 1.054203980480429, 1.015795900169476),
    @(1.02, 1.028964195235545, 1.039383125403378, 1.039201087819551, 1.050769629923122, 1.035593002171431, 1.03554617861446, 1.042922012888143, 1.042740642272621, 1.054267305337483, 1.015813611808298),
    @(1.02, 1.029491525669184, 1.039787809601182, 1.039657150258927, 1.051251977617608, 1.035677216548218, 1.035853270782017, 1.043211876779958, 1.043081679433988, 1.054636001029191, 1.015916704570851),
    @(1.02, 1.029822576310302, 1.040041839657165, 1.039943577166721, 1.051554863555001, 1.03572968475908, 1.036045954896317, 1.043393654337817, 1.043295731406119, 1.054867376304377, 1.0159813744924),
    @(1.02, 1.0300157693681, 1.040190076174885, 1.040110770987376, 1.051731647146688, 1.035760159604269, 1.036158363052161, 1.043499664914821, 1.043420630065675, 1.055002369862939, 1.016019096121428),
    @(1.02, 1.030081659578209, 1.040240631976126, 1.04016780104552, 1.051791945283244, 1.035770528840932, 1.036196694470767, 1.043535808820242, 1.043463224996491, 1.055048405277799, 1.016031958335242),
    @(1.02, 1.029787047682591, 1.040014577890149, 1.039912833238484, 1.051522354813957, 1.035724068751904, 1.036025279762841, 1.043374153094322, 1.043272760910357, 1.054842548162407, 1.015974435939795),
    @(1.02, 1.028830377752513, 1.039280423580727, 1.039085391892853, 1.050647250659232, 1.035571508248101, 1.0354682176928, 1.042848395899818, 1.042654085132325, 1.054173717047742, 1.015787435120709),
    @(1.02, 1.028229765508024, 1.038819431506633, 1.038566294168932, 1.050098091790673, 1.035474430390646, 1.035118149209562, 1.042517686704831, 1.04226552206604, 1.05375353515896, 1.015669870600676),
    @(1.02, 1.028548077043156, 1.039063754609705, 1.038841368261113, 1.050389110693685, 1.035526002632651, 1.035303709754078, 1.042693015137742, 1.042471466092266, 1.053976248835345, 1.015732192630106),
    @(1.02, 1.029803101241194, 1.040026896108352, 1.039926724707699, 1.051537043758581, 1.035726606782407, 1.036034621907879, 1.043382964920159, 1.043283140147429, 1.054853766819187, 1.015977571170773),
    @(1.02, 1.031263278797502, 1.041147114190665, 1.041191150140472, 1.052873673373429, 1.035954334737452, 1.036883536203409, 1.044182933292742, 1.044226832404816, 1.055873494295358, 1.016262347203186)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowNum = $r + 2
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $addr = $cols[$c] + $rowNum
        $ws.Range($addr).Value = $rowValues[$c]
    }
}
